$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column formatting: mark the defined columns (A:K) and the trailing
# catch-all column (L:XFD) as collapsed, matching the worksheet's <cols>
# definition. ---
$ws.Columns.Item(1).ShowDetail = $false
$ws.Columns.Item(2).ShowDetail = $false
$ws.Columns.Item(3).ShowDetail = $false
$ws.Columns.Item(4).ShowDetail = $false
$ws.Columns.Item(5).ShowDetail = $false
$ws.Columns.Item(6).ShowDetail = $false
$ws.Columns.Item(7).ShowDetail = $false
$ws.Columns.Item(8).ShowDetail = $false
$ws.Columns.Item(9).ShowDetail = $false
$ws.Columns.Item(10).ShowDetail = $false
$ws.Columns.Item(11).ShowDetail = $false
$ws.Range("L1").EntireColumn.ShowDetail = $false

# --- Data edit: rule "R10" (row 8) "From" value changes from 0 to 2 ---
$ws.Range("C8").Value = 2
